$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (number format, font, border, alignment) from the last
# existing date cell so the new column-A cells reuse the same style index
# instead of Excel minting a brand-new one.
$dateTemplate = $ws.Cells.Item(119, 1)

$data = @(
    @{ Row = 120; A = 45550; B = 601.5386009199999;  C = 160.589920235;   I = 235.210510047;  K = 289.28176403527;  N = 36.50686441472;  O = 0.022202408; Q = [double]"1.764E-06";  U = 252.2760237565012; Z = 223.210321692912 },
    @{ Row = 121; A = 45551; B = 592.1998596118999;  C = 159.174071968;   I = 235.389554357;  K = 289.087615200347; N = 37.72919246432;  O = 0.021412044; Q = [double]"1.728E-06";  U = 245.4957857954999; Z = 220.23501632359 },
    @{ Row = 122; A = 45552; B = 613.5627606118999;  C = 162.37186443;    I = 235.514885374;  K = 291.2232523845;   N = 40.52308514912;  O = 0.021913544; Q = [double]"1.7592E-06"; U = 252.5318817927654; Z = 218.508090829508 },
    @{ Row = 123; A = 45553; B = 546.6240842922;     C = 164.6564971625;  I = 240.528126054;  K = 290.446657044808; N = 41.62900100352;  O = 0.02244714;  Q = [double]"1.8072E-06"; U = 260.3355518988235; Z = 322.801671045952 }
)

foreach ($entry in $data) {
    $r = $entry.Row

    $cellA = $ws.Cells.Item($r, 1)
    $dateTemplate.Copy($cellA)
    $cellA.Value = $entry.A

    $ws.Cells.Item($r, 2).Value  = $entry.B   # B
    $ws.Cells.Item($r, 3).Value  = $entry.C   # C
    $ws.Cells.Item($r, 9).Value  = $entry.I   # I
    $ws.Cells.Item($r, 11).Value = $entry.K   # K
    $ws.Cells.Item($r, 14).Value = $entry.N   # N
    $ws.Cells.Item($r, 15).Value = $entry.O   # O
    $ws.Cells.Item($r, 17).Value = $entry.Q   # Q
    $ws.Cells.Item($r, 21).Value = $entry.U   # U
    $ws.Cells.Item($r, 26).Value = $entry.Z   # Z
}
